$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, font, borders, alignment) of the
# existing "2021" column (R) into the new "2022" column (S) so the new
# column visually matches the rest of the table.
$ws.Range("R4:R14").Copy()
$ws.Range("S4:S14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new "2022" data column.
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 4.9538761752705343
$ws.Range("S6").Value = 11.304954640614097
$ws.Range("S7").Value = 5.1593323216995444
$ws.Range("S8").Value = 13.687943262411348
$ws.Range("S9").Value = 10.22864019253911
$ws.Range("S10").Value = 9.1213700670141478
$ws.Range("S11").Value = 3.1335149863760217
$ws.Range("S12").Value = 2.872905173311127
$ws.Range("S13").Value = 3.527842284697861
$ws.Range("S14").Value = 5.0305321314335565

# Match the author's final selection in the sheet.
$ws.Range("T6").Select()
